$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("NewLoanInput")
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsTrans = $wb.Worksheets.Item("Transactions")

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. Style fix-ups: the workbook author removed the custom "0.00" number
#    format style (old index 12) and reassigned every cell that used it to
#    the general-purpose wrap/center style (old index 7), except for two
#    cells which moved to the "#,##0.00" style (old index 10).  We replicate
#    this purely through copy / paste-special (formats only) so the engine
#    reuses the existing style entries instead of fabricating new ones.
# ---------------------------------------------------------------------------

# Grab a stable "style 10" source (Summary!F2) BEFORE its own style changes.
$wsSummary.Range("F2").Copy()
$wsSummary.Range("E2").PasteSpecial($xlPasteFormats)

# Grab a stable "style 7" source (Summary!C2) and fan it out everywhere a
# style-12 (or, in a couple of spots, a now-unwanted style-10) cell needs to
# become plain style 7.
$wsSummary.Range("C2").Copy()
$wsSummary.Range("B2").PasteSpecial($xlPasteFormats)
$wsSummary.Range("A3").PasteSpecial($xlPasteFormats)
$wsSummary.Range("B3").PasteSpecial($xlPasteFormats)
$wsSummary.Range("E3").PasteSpecial($xlPasteFormats)
$wsSummary.Range("F3").PasteSpecial($xlPasteFormats)

# Repayment schedule: F/H/P columns (rows 3, 5-15) + P4 move from style 12
# to style 7; K column (rows 3-15) plus L3/M3 move from style 10 to style 7.
foreach ($r in 3,5,6,7,8,9,10,11,12,13,14,15) {
  $wsSchedule.Range("F$r").PasteSpecial($xlPasteFormats)
  $wsSchedule.Range("H$r").PasteSpecial($xlPasteFormats)
  $wsSchedule.Range("K$r").PasteSpecial($xlPasteFormats)
  $wsSchedule.Range("P$r").PasteSpecial($xlPasteFormats)
}
$wsSchedule.Range("P4").PasteSpecial($xlPasteFormats)
$wsSchedule.Range("L3").PasteSpecial($xlPasteFormats)
$wsSchedule.Range("M3").PasteSpecial($xlPasteFormats)

# Transactions: E3/F3/G3 move from style 12 to style 7.
$wsTrans.Range("E3").PasteSpecial($xlPasteFormats)
$wsTrans.Range("F3").PasteSpecial($xlPasteFormats)
$wsTrans.Range("G3").PasteSpecial($xlPasteFormats)

# Transactions: J2/J3 pick up the "#,##0.00" style 10 (borrowed from the
# Repayment schedule sheet, where that style is still very much in use).
$wsSchedule.Range("G5").Copy()
$wsTrans.Range("J2").PasteSpecial($xlPasteFormats)
$wsTrans.Range("J3").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Summary sheet: drop the stray empty G2 cell (column G disappears from
#    the sheet's used range entirely) and refresh the values that moved.
# ---------------------------------------------------------------------------
$wsSummary.Range("G2").Clear()

$wsSummary.Range("A3").Value = 415.28
$wsSummary.Range("B3").Value = 16.670000000000002
$wsSummary.Range("E3").Value = 398.61
$wsSummary.Range("F3").Value = 25

# ---------------------------------------------------------------------------
# 3. Repayment schedule: drop the now-unused "O" column (Over Due helper)
#    across every data row, and refresh the recomputed figures.
# ---------------------------------------------------------------------------
foreach ($r in 2,3,4,5,6,7,8,9,10,11,12,13,14,15) {
  $wsSchedule.Range("O$r").Clear()
}

$wsSchedule.Range("G3").Value = 4166.67
$wsSchedule.Range("H3").Value = 16.670000000000002
$wsSchedule.Range("K3").Value = 850
$wsSchedule.Range("L3").Value = 850
$wsSchedule.Range("M3").Value = 0
$wsSchedule.Range("N3").Value = 850

$wsSchedule.Range("K4").Value = 0

$wsSchedule.Range("G5").Value = 8333.34
$wsSchedule.Range("H5").Value = 25
$wsSchedule.Range("K5").Value = 858.33
$wsSchedule.Range("P5").Value = 858.33

$wsSchedule.Range("G6").Value = 7500.01
$wsSchedule.Range("H6").Value = 61.11
$wsSchedule.Range("K6").Value = 894.44
$wsSchedule.Range("P6").Value = 894.44

$wsSchedule.Range("G7").Value = 6666.68
$wsSchedule.Range("H7").Value = 62.5
$wsSchedule.Range("K7").Value = 895.83
$wsSchedule.Range("P7").Value = 895.83

$wsSchedule.Range("G8").Value = 5833.35
$wsSchedule.Range("H8").Value = 55.56
$wsSchedule.Range("K8").Value = 888.89
$wsSchedule.Range("P8").Value = 888.89

$wsSchedule.Range("G9").Value = 5000.0200000000004
$wsSchedule.Range("H9").Value = 48.61
$wsSchedule.Range("K9").Value = 881.94
$wsSchedule.Range("P9").Value = 881.94

$wsSchedule.Range("G10").Value = 4166.6899999999996
$wsSchedule.Range("H10").Value = 41.67
$wsSchedule.Range("K10").Value = 875
$wsSchedule.Range("P10").Value = 875

$wsSchedule.Range("G11").Value = 3333.36
$wsSchedule.Range("H11").Value = 34.72
$wsSchedule.Range("K11").Value = 868.05
$wsSchedule.Range("P11").Value = 868.05

$wsSchedule.Range("G12").Value = 2500.0300000000002
$wsSchedule.Range("H12").Value = 27.78
$wsSchedule.Range("K12").Value = 861.11
$wsSchedule.Range("P12").Value = 861.11

$wsSchedule.Range("G13").Value = 1666.7
$wsSchedule.Range("H13").Value = 20.83
$wsSchedule.Range("K13").Value = 854.16
$wsSchedule.Range("P13").Value = 854.16

$wsSchedule.Range("G14").Value = 833.37
$wsSchedule.Range("H14").Value = 13.89

$wsSchedule.Range("F15").Value = 833.37
$wsSchedule.Range("H15").Value = 6.94
$wsSchedule.Range("K15").Value = 840.31
$wsSchedule.Range("P15").Value = 840.31

# ---------------------------------------------------------------------------
# 4. Transactions sheet: renumbered transaction IDs and recomputed balances.
# ---------------------------------------------------------------------------
$wsTrans.Range("A2").Value = 108
$wsTrans.Range("J2").Value = 9166.67

$wsTrans.Range("A3").Value = 106
$wsTrans.Range("E3").Value = 850
$wsTrans.Range("J3").Value = 4166.67

$wsTrans.Range("A4").Value = 104

# ---------------------------------------------------------------------------
# 5. View state: reselect ranges on each sheet and make "Transactions" the
#    active tab (was "Repayment schedule").
# ---------------------------------------------------------------------------
$wsSummary.Activate()
$wsSummary.Range("A7:XFD13").Select()

$wsSchedule.Activate()
$wsSchedule.Range("A16:XFD16").Select()

$wsTrans.Activate()
$wsTrans.Range("A2:XFD6").Select()

Write-Host "Edit complete"
